# ---------------------------------------------------------------------------
# Applies the "finished second corr calcuations" commit:
#   1. Renames "Correlation" sheet -> "Box Office $ Correlation"
#   2. Adds a new sheet "Rank + In Theaters Corr" at the end containing a
#      copy of the Weekly Data (Published/Week-End/Rank) columns plus a new
#      "In Theaters?" flag column and two CORREL() summaries.
#   3. Widens the selected range on the "NYT Books" sheet from B1:D77 to
#      B1:E77.
#   4. Makes the new sheet the active tab.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- data for the new sheet -------------------------------------------------
$Bvals = @(41944,41951,41958,41965,41972,41979,41986,41993,42000,42007,42014,42021,42028,42035,42042,42049,42056,42063,42070,42077,42084,42091,42098,42105,42112,42119,42126,42133,42140,42147,42154,42161,42168,42175,42182,42189,42196,42203,42210,42217,42224,42231,42238,42245,42252,42259,42266,42273,42280,42287,42294,42301,42308,42315,42322,42329,42336,42343,42350,42357,42364,42371,42378,42385,42392,42399,42406,42413,42420,42427,42434,42441,42448,42455,42462,42469)
$Cvals = @(8,7,8,7,5,4,3,4,3,2,4,5,6,6,7,8,7,7,9,7,8,4,5,7,4,5,6,6,9,7,5,7,1,2,2,2,2,2,2,2,2,2,1,1,1,1,1,1,1,1,1,1,1,2,2,1,1,1,1,1,1,1,2,2,2,2,4,3,4,5,5,6,5,7,11,13)
$Dvals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0)

# --- 1. rename the existing Correlation sheet -------------------------------
$corrSheet = $wb.Worksheets.Item("Correlation")
$corrSheet.Name = "Box Office `$ Correlation"

# --- 2. add the new sheet at the end of the workbook ------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Rank + In Theaters Corr"

# column widths (ColumnWidth is offset from the raw OOXML "width" by 5/6)
$newSheet.Columns.Item(2).ColumnWidth = 14.5 - 0.8333333333333334
$newSheet.Columns.Item(3).ColumnWidth = 14.5 - 0.8333333333333334
$newSheet.Columns.Item(4).ColumnWidth = 22.83203125 - 0.8333333333333334
$newSheet.Columns.Item(6).ColumnWidth = 33.5 - 0.8333333333333334

# header row
$newSheet.Cells.Item(1,1).Value2 = "Week Start"
$newSheet.Cells.Item(1,2).Value2 = "Week End"
$newSheet.Cells.Item(1,3).Value2 = "Rank"
$newSheet.Cells.Item(1,4).Value2 = "In Theaters? (1 = Y, 0 = N)"
$newSheet.Cells.Item(1,6).Value2 = "Correlation:"

$newSheet.Cells.Item(1,2).NumberFormat = "yyyy\-mm\-dd;@"

# data rows 2-77: Week Start (A, formula), Week End (B), Rank (C), In Theaters (D)
for ($i = 0; $i -lt $Bvals.Length; $i++) {
  $r = $i + 2
  $newSheet.Cells.Item($r,1).Formula = "=B$r-6"
  $newSheet.Cells.Item($r,1).NumberFormat = "yyyy\-mm\-dd;@"
  $newSheet.Cells.Item($r,2).Value2 = $Bvals[$i]
  $newSheet.Cells.Item($r,2).NumberFormat = "yyyy\-mm\-dd;@"
  $newSheet.Cells.Item($r,3).Value2 = $Cvals[$i]
  $newSheet.Cells.Item($r,4).Value2 = $Dvals[$i]
}

# correlation formulas
$newSheet.Cells.Item(2,6).Formula = "=CORREL(C2:C77,D2:D77)"
$newSheet.Cells.Item(4,6).Value2 = "Correlation once first move comes out:"
$newSheet.Cells.Item(5,6).Formula = "=CORREL(C50:C77,D50:D77)"

# --- 3. widen the selection on the NYT Books sheet --------------------------
$booksSheet = $wb.Worksheets.Item("NYT Books")
$booksSheet.Range("B1:E77").Select()

# --- 4. make the new sheet the active tab + set its selection ---------------
$newSheet.Activate()
$newSheet.Range("F6").Select()

Write-Output "done"
